# CalSim3GroundWaterDataExtractionInitFile_v1.xlsx
# "updated init files up to scenario 56"
#
# The Upper-Left/Lower-Right "block" helper cells in column D referenced
# row 28 of the various scenario-listing blocks; the edit moves those
# references down to row 31 (reflecting the listing sheet now running
# through "scenario 56"). Two other lookup-range helper cells (D15/D22)
# also move to reflect the grown source ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D "Lower Right Cell" reference updates -----------------------
$ws.Range("D5").Value  = "A31"
$ws.Range("D6").Value  = "B31"
$ws.Range("D7").Value  = "C31"
$ws.Range("D8").Value  = "G31"
$ws.Range("D9").Value  = "H31"
$ws.Range("D10").Value = "I31"
$ws.Range("D11").Value = "J31"

$ws.Range("D15").Value = "E258"
$ws.Range("D22").Value = "O440"

# --- New (blank, bold+italic styled) helper cells in column E -----------
# Mirrors the existing column-D styling (bold+italic) so the new cells
# pick up the same cell style already used for D5:D11.
foreach ($r in 5..11) {
    $cell = $ws.Range("E$r")
    $cell.Font.Bold = $true
    $cell.Font.Italic = $true
}

# --- Selection moves to reflect the newly highlighted D5:D11 block ------
[void]$ws.Range("D5:D11").Select()
